$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.774.41'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.529.15'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.37%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '624.94'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.24'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("E7").Value = '  -0.84%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.526.50'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.200'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.04'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -5.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.583'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.49'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000280'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.096.15'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.42%  '

$ws.Range("E16").Value = '  -0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '607.20'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.524.26'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.825.08'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.884'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.08'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.36%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '98.22'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.53%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.64'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.77'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.10%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("E28").Value = '  -2.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.80'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.06'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.58%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.13'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.31'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.61%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.80'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.42%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '626.27'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -14.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0997'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.81'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("E38").Value = '  -9.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0474'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.87'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.83%  '

$ws.Range("E41").Value = '  +0.27%  '

$ws.Range("E42").Value = '  +1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.357.96'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.00%  '

$ws.Range("E44").Value = '  +2.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.97'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("E46").Value = '  -3.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.05'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.92%  '

$ws.Range("E48").Value = '  -3.32%  '

$ws.Range("E49").Value = '  -0.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.20'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.14%  '

$ws.Range("E51").Value = '  +6.21%  '

